$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 41
$ws.Range("B12").Value = 8
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 11
$ws.Range("B15").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 1
$ws.Range("B24").Value = 1
$ws.Range("B29").Value = 14
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 31
$ws.Range("B32").Value = 80
$ws.Range("B33").Value = 122
$ws.Range("B35").Value = 6
$ws.Range("B36").Value = 0
